$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.131.90'
$ws.Range('E2').Value = '  -0.50%  '
$ws.Range('D3').Value = '1.901.72'
$ws.Range('E3').Value = '  +0.04%  '
$ws.Range('E4').Value = '  -0.44%  '
$ws.Range('D5').Value = '''253.28'
$ws.Range('E5').Value = '  +3.13%  '
$ws.Range('D6').Value = '''0.694'
$ws.Range('E6').Value = '  +0.11%  '
$ws.Range('E7').Value = '  -0.38%  '
$ws.Range('D8').Value = '''41.24'
$ws.Range('E8').Value = '  +0.94%  '
$ws.Range('D9').Value = '''0.358'
$ws.Range('E9').Value = '  +3.08%  '
$ws.Range('D10').Value = '''52.81'
$ws.Range('E10').Value = '  -0.62%  '
$ws.Range('D11').Value = '''0.0750'
$ws.Range('E11').Value = '  +3.91%  '
$ws.Range('D12').Value = '''0.0982'
$ws.Range('E12').Value = '  -1.40%  '
$ws.Range('D13').Value = '''13.05'
$ws.Range('E13').Value = '  +5.63%  '
$ws.Range('E14').Value = '  +0.00%  '
$ws.Range('E15').Value = '  +3.99%  '
$ws.Range('E16').Value = '  +3.29%  '
$ws.Range('D17').Value = '1.892.69'
$ws.Range('E17').Value = '  -0.52%  '
$ws.Range('D18').Value = '35.131.06'
$ws.Range('D19').Value = '''73.51'
$ws.Range('E19').Value = '  +1.71%  '
$ws.Range('D20').Value = '0.0₃0832'
$ws.Range('E20').Value = '  +1.65%  '
$ws.Range('D21').Value = '''242.84'
$ws.Range('E21').Value = '  +0.84%  '
$ws.Range('D22').Value = '''12.96'
$ws.Range('E22').Value = '  +2.68%  '
$ws.Range('D23').Value = '''5.04'
$ws.Range('E23').Value = '  +5.19%  '
$ws.Range('E24').Value = '  -0.40%  '
$ws.Range('D25').Value = '''2.44'
$ws.Range('E25').Value = '  +5.56%  '
$ws.Range('E26').Value = '  +0.78%  '
$ws.Range('D27').Value = '''167.00'
$ws.Range('E27').Value = '  -0.81%  '
$ws.Range('E28').Value = '  -0.34%  '
$ws.Range('E29').Value = '  +0.97%  '
$ws.Range('D30').Value = '''0.129'
$ws.Range('E30').Value = '  -0.84%  '
$ws.Range('D31').Value = '4.128.06'
$ws.Range('E31').Value = '  -0.28%  '
$ws.Range('E32').Value = '  +13.89%  '
$ws.Range('D33').Value = '''0.0609'
$ws.Range('E33').Value = '  +6.97%  '
$ws.Range('E34').Value = '  +4.12%  '
$ws.Range('D35').Value = '''1.58'
$ws.Range('E35').Value = '  +8.71%  '
$ws.Range('E36').Value = '  +2.76%  '
$ws.Range('E37').Value = '  -0.44%  '
$ws.Range('D38').Value = '''0.854'
$ws.Range('E38').Value = '  -7.59%  '
$ws.Range('E39').Value = '  -1.09%  '
$ws.Range('D40').Value = '''103.14'
$ws.Range('E40').Value = '  +14.67%  '
$ws.Range('D41').Value = '''17.35'
$ws.Range('E41').Value = '  +8.29%  '
$ws.Range('E42').Value = '  +2.76%  '
$ws.Range('E43').Value = '  +0.92%  '
$ws.Range('D44').Value = '''0.0650'
$ws.Range('E44').Value = '  -1.61%  '
$ws.Range('E45').Value = '  +0.03%  '
$ws.Range('D46').Value = '1.319.48'
$ws.Range('E46').Value = '  -2.30%  '
$ws.Range('E47').Value = '  +0.02%  '
$ws.Range('E48').Value = '  -1.54%  '
$ws.Range('E49').Value = '  +1.09%  '
$ws.Range('D50').Value = '''11.95'
$ws.Range('E50').Value = '  -5.88%  '
$ws.Range('D51').Value = '''0.0745'
$ws.Range('E51').Value = '  +5.73%  '
